# Trade #51 closed at 2026-02-17 08:39:46 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status rollups and appends the new
# MarketMaking trade record (#51) to both the "All Trades" and
# "MarketMaking" detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet roll-up figures
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.46   # Current Capital
$summary.Range("B4").Value = -0.54     # Total P&L $
$summary.Range("B5").Value = -0.21     # Total P&L %
$summary.Range("B6").Value = 51        # Total Trades
$summary.Range("B8").Value = 23        # Losing Trades
$summary.Range("B9").Value = 35.29     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.45999999999999   # Capital
$status.Range("D4").Value = 51                  # Trades
$status.Range("E4").Value = -0.54               # P&L $
$status.Range("F4").Value = -0.54               # P&L %
$status.Range("G4").Value = 35.29               # Win Rate %

# ---------------------------------------------------------------------
# 3) Append trade #51 to a detail sheet (All Trades / MarketMaking share
#    the same row layout); helper writes the literal text columns so
#    the date-like string isn't auto-coerced into a date serial, while
#    keeping default (unstyled) cells like the source rows.
# ---------------------------------------------------------------------
function Add-TradeRow51($ws) {
    $row = 52

    $ws.Cells.Item($row, 1).Value = 51          # Trade #

    $dateCell = $ws.Cells.Item($row, 2)         # Date (keep as text)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 3).Value = "08:39:40"  # Time
    $ws.Cells.Item($row, 4).Value = "MarketMaking"  # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"          # Side
    $ws.Cells.Item($row, 6).Value = 0.09            # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.04            # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"        # Status
    $ws.Cells.Item($row, 9).Value = -55.5556        # P&L %
    $ws.Cells.Item($row, 10).Value = -0.05          # P&L $
    $ws.Cells.Item($row, 11).Value = 99.45999999999999  # Capital After
    $ws.Cells.Item($row, 12).Value = 0              # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0              # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6            # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"   # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.13           # Duration (min)
}

Add-TradeRow51 $wb.Worksheets.Item("All Trades")
Add-TradeRow51 $wb.Worksheets.Item("MarketMaking")

Write-Output "Trade #51 appended; summary + strategy status updated."
